# Refresh the team-specific transition-probability matrix on Sheet1 with the
# recomputed values after adding more simulated games (per commit message:
# "added more games, sped up simulate game logic, and drafted optimization
# logic"). Each changed cell below is set to the exact recomputed
# probability from the updated simulation run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2321428571428572
$ws.Range("C2").Value = 0.5029761904761905
$ws.Range("J2").Value = 0.002976190476190476
$ws.Range("P2").Value = 0.1517857142857143
$ws.Range("S2").Value = 0.1101190476190476
$ws.Range("C3").Value = 0.01764705882352941
$ws.Range("J3").Value = 0.02352941176470588
$ws.Range("P3").Value = 0.7647058823529411
$ws.Range("S3").Value = 0.1941176470588235
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.04761904761904762
$ws.Range("D6").Value = 0.01058201058201058
$ws.Range("F6").Value = 0.06349206349206349
$ws.Range("J6").Value = 0.2804232804232804
$ws.Range("O6").Value = 0.02645502645502645
$ws.Range("Q6").Value = 0.1481481481481481
$ws.Range("R6").Value = 0.05291005291005291
$ws.Range("S6").Value = 0.3703703703703703
$ws.Range("B7").Value = 0.1813186813186813
$ws.Range("D7").Value = 0.02197802197802198
$ws.Range("F7").Value = 0.02747252747252747
$ws.Range("J7").Value = 0.1043956043956044
$ws.Range("O7").Value = 0.02197802197802198
$ws.Range("Q7").Value = 0.1923076923076923
$ws.Range("R7").Value = 0.07142857142857142
$ws.Range("S7").Value = 0.3791208791208791
$ws.Range("B8").Value = 0.09941520467836257
$ws.Range("D8").Value = 0.01364522417153996
$ws.Range("E8").Value = 0.003898635477582846
$ws.Range("F8").Value = 0.04093567251461988
$ws.Range("J8").Value = 0.1091617933723197
$ws.Range("O8").Value = 0.01754385964912281
$ws.Range("Q8").Value = 0.1715399610136452
$ws.Range("R8").Value = 0.07407407407407407
$ws.Range("S8").Value = 0.4697855750487329
$ws.Range("B9").Value = 0.09547738693467336
$ws.Range("D9").Value = 0.02010050251256281
$ws.Range("F9").Value = 0.05527638190954774
$ws.Range("J9").Value = 0.1306532663316583
$ws.Range("Q9").Value = 0.1658291457286432
$ws.Range("R9").Value = 0.1055276381909548
$ws.Range("S9").Value = 0.4271356783919598
$ws.Range("B10").Value = 0.1243291592128801
$ws.Range("D10").Value = 0.01788908765652952
$ws.Range("F10").Value = 0.07155635062611806
$ws.Range("J10").Value = 0.10912343470483
$ws.Range("O10").Value = 0.009838998211091235
$ws.Range("Q10").Value = 0.1994633273703041
$ws.Range("R10").Value = 0.06618962432915922
$ws.Range("S10").Value = 0.4016100178890877
$ws.Range("G11").Value = 0.1638225255972696
$ws.Range("J11").Value = 0.10580204778157
$ws.Range("K11").Value = 0.2320819112627986
$ws.Range("L11").Value = 0.4778156996587031
$ws.Range("S11").Value = 0.0204778156996587
$ws.Range("G12").Value = 0.8111888111888111
$ws.Range("J12").Value = 0.1608391608391608
$ws.Range("K12").Value = 0.006993006993006993
$ws.Range("L12").Value = 0.01398601398601399
$ws.Range("S12").Value = 0.006993006993006993
$ws.Range("G13").Value = 0.5625
$ws.Range("J13").Value = 0.3125
$ws.Range("S13").Value = 0.125
$ws.Range("F15").Value = 0.01092896174863388
$ws.Range("H15").Value = 0.1748633879781421
$ws.Range("I15").Value = 0.1038251366120219
$ws.Range("J15").Value = 0.360655737704918
$ws.Range("K15").Value = 0.04918032786885246
$ws.Range("M15").Value = 0.01092896174863388
$ws.Range("O15").Value = 0.0273224043715847
$ws.Range("S15").Value = 0.2622950819672131
$ws.Range("F16").Value = 0.02450980392156863
$ws.Range("H16").Value = 0.1813725490196078
$ws.Range("I16").Value = 0.1225490196078431
$ws.Range("J16").Value = 0.3529411764705883
$ws.Range("K16").Value = 0.1127450980392157
$ws.Range("M16").Value = 0.03431372549019608
$ws.Range("O16").Value = 0.03431372549019608
$ws.Range("S16").Value = 0.1372549019607843
$ws.Range("F17").Value = 0.01237623762376238
$ws.Range("H17").Value = 0.2376237623762376
$ws.Range("I17").Value = 0.1014851485148515
$ws.Range("J17").Value = 0.3663366336633663
$ws.Range("K17").Value = 0.07425742574257425
$ws.Range("M17").Value = 0.01485148514851485
$ws.Range("O17").Value = 0.04950495049504951
$ws.Range("S17").Value = 0.1435643564356436
$ws.Range("F18").Value = 0.02564102564102564
$ws.Range("H18").Value = 0.2564102564102564
$ws.Range("I18").Value = 0.07692307692307693
$ws.Range("J18").Value = 0.3782051282051282
$ws.Range("K18").Value = 0.07051282051282051
$ws.Range("M18").Value = 0.01923076923076923
$ws.Range("O18").Value = 0.05128205128205128
$ws.Range("S18").Value = 0.1217948717948718
$ws.Range("F19").Value = 0.01601830663615561
$ws.Range("H19").Value = 0.2349351639969489
$ws.Range("I19").Value = 0.08161708619374523
$ws.Range("J19").Value = 0.3356216628527841
$ws.Range("K19").Value = 0.1113653699466056
$ws.Range("M19").Value = 0.02364607170099161
$ws.Range("O19").Value = 0.06636155606407322
$ws.Range("S19").Value = 0.1304347826086956
